$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 values (memmove_dma / my_memmove_dma row under FRDM and BBB tables)
$ws.Range("C7").Value = 187
$ws.Range("D7").Value = 1267
$ws.Range("E7").Value = 12067
$ws.Range("F7").Value = 60067

$ws.Range("I7").Value = 152
$ws.Range("J7").Value = 1052
$ws.Range("K7").Value = 10052
$ws.Range("L7").Value = 50052

# Row 12 header label cells (unit label "uS" for the uS table)
$ws.Range("B12").Value = "uS"
$ws.Range("H12").Value = "uS"

# Row 13 values
$ws.Range("C13").Value = 25
$ws.Range("D13").Value = 6
$ws.Range("E13").Value = 8
$ws.Range("F13").Value = 25

$ws.Range("I13").Value = 18
$ws.Range("J13").Value = 5
$ws.Range("K13").Value = 5
$ws.Range("L13").Value = 13

# Row 14 values
$ws.Range("C14").Value = 7
$ws.Range("D14").Value = 17
$ws.Range("E14").Value = 128
$ws.Range("F14").Value = 623

$ws.Range("I14").Value = 6
$ws.Range("J14").Value = 13
$ws.Range("K14").Value = 94
$ws.Range("L14").Value = 454

# Row 15 values
$ws.Range("C15").Value = 5
$ws.Range("D15").Value = 7
$ws.Range("E15").Value = 28
$ws.Range("F15").Value = 122

$ws.Range("I15").Value = 4
$ws.Range("J15").Value = 6
$ws.Range("K15").Value = 20
$ws.Range("L15").Value = 88

# Sheet view changes: remove topLeftCell freeze/scroll position, update selection
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Application.ActiveWindow.ScrollColumn = 1
$null = $ws.Range("F8").Select()
